$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 188.83333
$ws.Range("I33").Value = 188.09091
$ws.Range("K33").Value = 188.09091
$ws.Range("M33").Value = 40.90908999999999
$ws.Range("H41").Value = 652.82355
$ws.Range("I41").Value = 329.66666
$ws.Range("K41").Value = 329.66666
$ws.Range("M41").Value = 110.33334
$ws.Range("H92").Value = 2340.2856
$ws.Range("I92").Value = 2205.9092
$ws.Range("K92").Value = 2205.9092
$ws.Range("M92").Value = -957.9092000000001
$ws.Range("H98").Value = 644.64703
$ws.Range("I98").Value = 644.64703
$ws.Range("K98").Value = 644.64703
$ws.Range("M98").Value = 853.35297
$ws.Range("H113").Value = 3424.5833
$ws.Range("I113").Value = 2391.923
$ws.Range("K113").Value = 2391.923
$ws.Range("M113").Value = 862.0770000000002
$ws.Range("H122").Value = 644.64703
$ws.Range("I122").Value = 644.64703
$ws.Range("K122").Value = 1933.94109
$ws.Range("M122").Value = 516.0589100000002
$ws.Range("H132").Value = 1448.0476
$ws.Range("I132").Value = 1318.3529
$ws.Range("K132").Value = 3955.0587
$ws.Range("M132").Value = -1425.0587
$ws.Range("H133").Value = 76107.164
$ws.Range("J133").Value = 76107.164
$ws.Range("L133").Value = 76107.164
$ws.Range("N133").Value = -86227.164
$ws.Range("H134").Value = 74143.164
$ws.Range("J134").Value = 74143.164
$ws.Range("L134").Value = 74143.164
$ws.Range("N134").Value = -84283.164
$ws.Range("H135").Value = 2010.6666
$ws.Range("I135").Value = 2219.2
$ws.Range("K135").Value = 19972.8
$ws.Range("M135").Value = -17437.8
$ws.Range("H136").Value = 77987
$ws.Range("J136").Value = 77987
$ws.Range("L136").Value = 77987
$ws.Range("N136").Value = -88187
$ws.Range("H139").Value = 74223.125
$ws.Range("J139").Value = 74223.125
$ws.Range("L139").Value = 74223.125
$ws.Range("N139").Value = -84503.125
$ws.Range("H140").Value = 50556.715
$ws.Range("J140").Value = 50531.332
$ws.Range("L140").Value = 50531.332
$ws.Range("N140").Value = -60891.332

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 790.9091
$ws.Range("I97").Value = 810
$ws.Range("K97").Value = 810
$ws.Range("M97").Value = -314
$ws.Range("H139").Value = 95999.5
$ws.Range("J139").Value = 95999.5
$ws.Range("L139").Value = 95999.5
$ws.Range("N139").Value = -106279.5

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2310.3333
$ws.Range("I86").Value = 2157.0908
$ws.Range("J86").Value = 2731.75
$ws.Range("K86").Value = 2157.0908
$ws.Range("L86").Value = 2731.75
$ws.Range("M86").Value = -1034.0908
$ws.Range("N86").Value = -4977.75
$ws.Range("H89").Value = 2310.3333
$ws.Range("I89").Value = 2157.0908
$ws.Range("J89").Value = 2731.75
$ws.Range("K89").Value = 10785.454
$ws.Range("L89").Value = 13658.75
$ws.Range("M89").Value = -5169.454
$ws.Range("N89").Value = -24890.75
$ws.Range("H132").Value = 28704.092
$ws.Range("J132").Value = 28704.092
$ws.Range("L132").Value = 28704.092
$ws.Range("N132").Value = -38824.092
$ws.Range("H134").Value = 2774.2703
$ws.Range("I134").Value = 1614.125
$ws.Range("K134").Value = 4842.375
$ws.Range("M134").Value = -2307.375
$ws.Range("H138").Value = 72152.25
$ws.Range("J138").Value = 72152.25
$ws.Range("L138").Value = 72152.25
$ws.Range("N138").Value = -82432.25
$ws.Range("H140").Value = 65261
$ws.Range("J140").Value = 74796.8
$ws.Range("L140").Value = 74796.8
$ws.Range("N140").Value = -85156.8

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1551.7778
$ws.Range("I16").Value = 1110.091
$ws.Range("K16").Value = 1110.091
$ws.Range("M16").Value = -823.0909999999999
$ws.Range("H31").Value = 2817.2307
$ws.Range("I31").Value = 2147
$ws.Range("J31").Value = 4157.6924
$ws.Range("K31").Value = 2147
$ws.Range("L31").Value = 4157.6924
$ws.Range("M31").Value = -1852
$ws.Range("N31").Value = -4747.6924
$ws.Range("H34").Value = 2817.2307
$ws.Range("I34").Value = 2147
$ws.Range("J34").Value = 4157.6924
$ws.Range("K34").Value = 2147
$ws.Range("L34").Value = 4157.6924
$ws.Range("M34").Value = -1945
$ws.Range("N34").Value = -4561.6924
$ws.Range("H58").Value = 1644.8823
$ws.Range("I58").Value = 1444.2916
$ws.Range("J58").Value = 2126.3
$ws.Range("K58").Value = 1444.2916
$ws.Range("L58").Value = 2126.3
$ws.Range("M58").Value = -1241.2916
$ws.Range("N58").Value = -2532.3
$ws.Range("H113").Value = 1551.7778
$ws.Range("I113").Value = 1110.091
$ws.Range("K113").Value = 1110.091
$ws.Range("M113").Value = 1059.909
$ws.Range("H122").Value = 3035.48
$ws.Range("J122").Value = 2899.5293
$ws.Range("L122").Value = 8698.5879
$ws.Range("N122").Value = -13598.5879
$ws.Range("H132").Value = 1671986.5
$ws.Range("I132").Value = 2067907.1
$ws.Range("K132").Value = 6203721.300000001
$ws.Range("M132").Value = -6201191.300000001
$ws.Range("H134").Value = 3020574.8
$ws.Range("I134").Value = 3970974.8
$ws.Range("J134").Value = 169374.17
$ws.Range("K134").Value = 11912924.4
$ws.Range("L134").Value = 508122.51
$ws.Range("M134").Value = -11910389.4
$ws.Range("N134").Value = -513192.51
$ws.Range("H136").Value = 1644.8823
$ws.Range("I136").Value = 1444.2916
$ws.Range("J136").Value = 2126.3
$ws.Range("K136").Value = 4332.8748
$ws.Range("L136").Value = 6378.900000000001
$ws.Range("M136").Value = -1782.8748
$ws.Range("N136").Value = -11478.9

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 788.56525
$ws.Range("I5").Value = 765.9167
$ws.Range("J5").Value = 813.2727
$ws.Range("K5").Value = 2297.7501
$ws.Range("L5").Value = 2439.8181
$ws.Range("M5").Value = -2185.7501
$ws.Range("N5").Value = -2663.8181
$ws.Range("H103").Value = 395.9
$ws.Range("I103").Value = 395.9
$ws.Range("K103").Value = 1187.7
$ws.Range("M103").Value = -308.6999999999998
$ws.Range("H117").Value = 820.3333
$ws.Range("I117").Value = 1029
$ws.Range("J117").Value = 716
$ws.Range("K117").Value = 3087
$ws.Range("L117").Value = 2148
$ws.Range("M117").Value = 355
$ws.Range("N117").Value = -9032
$ws.Range("H135").Value = 788.56525
$ws.Range("I135").Value = 765.9167
$ws.Range("J135").Value = 813.2727
$ws.Range("K135").Value = 6893.2503
$ws.Range("L135").Value = 7319.454299999999
$ws.Range("M135").Value = -4358.2503
$ws.Range("N135").Value = -12389.4543

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 23167.834
$ws.Range("I21").Value = 19500
$ws.Range("J21").Value = 25001.75
$ws.Range("K21").Value = 19500
$ws.Range("L21").Value = 25001.75
$ws.Range("M21").Value = -19327
$ws.Range("N21").Value = -25347.75
$ws.Range("H30").Value = 23167.834
$ws.Range("I30").Value = 19500
$ws.Range("J30").Value = 25001.75
$ws.Range("K30").Value = 19500
$ws.Range("L30").Value = 25001.75
$ws.Range("M30").Value = -19395
$ws.Range("N30").Value = -25211.75
$ws.Range("H97").Value = 1253.091
$ws.Range("I97").Value = 362.4
$ws.Range("K97").Value = 362.4
$ws.Range("M97").Value = 133.6
$ws.Range("H102").Value = 1270.0834
$ws.Range("I102").Value = 1264.3
$ws.Range("K102").Value = 1264.3
$ws.Range("M102").Value = 357.7
$ws.Range("H107").Value = 796.7895
$ws.Range("I107").Value = 818
$ws.Range("K107").Value = 818
$ws.Range("M107").Value = 1102
$ws.Range("H109").Value = 27106.445
$ws.Range("J109").Value = 27106.445
$ws.Range("L109").Value = 27106.445
$ws.Range("N109").Value = -29186.445
$ws.Range("H113").Value = 3551.5
$ws.Range("J113").Value = 3800.4
$ws.Range("L113").Value = 3800.4
$ws.Range("N113").Value = -8140.4
$ws.Range("H135").Value = 90000
$ws.Range("J135").Value = 90000
$ws.Range("L135").Value = 90000
$ws.Range("N135").Value = -100140
$ws.Range("H140").Value = 94015
$ws.Range("J140").Value = 94334.55
$ws.Range("L140").Value = 94334.55
$ws.Range("N140").Value = -104694.55

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 75003840
$ws.Range("I122").Value = 71432510
$ws.Range("K122").Value = 214297530
$ws.Range("M122").Value = -214295080
$ws.Range("H136").Value = 2888.5557

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 2999.75
$ws.Range("I20").Value = 2999
$ws.Range("K20").Value = 2999
$ws.Range("M20").Value = -2759
$ws.Range("H122").Value = 1311.68
$ws.Range("I122").Value = 863.0526
$ws.Range("J122").Value = 2732.3333
$ws.Range("K122").Value = 2589.1578
$ws.Range("L122").Value = 8196.999899999999
$ws.Range("M122").Value = -139.1578
$ws.Range("N122").Value = -13096.9999
